$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price (D) and Volume (E) columns hold plain text values (e.g. "96.639.31",
# "  -0.86%  "). Force each touched cell to Text format before assigning so that
# numeric-looking strings are not auto-converted into actual numbers by Excel,
# keeping them as text just like the untouched cells in the same columns.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "96.642.04"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.88%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.715.15"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +3.17%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.18"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.31%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +9.47%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "656.87"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.51%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.04%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.24%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.714.10"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +3.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "45.27"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.04%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.58%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +5.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.406.54"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +3.13%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000270"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "96.610.47"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.68%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +17.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.696.98"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.69%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +4.71%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.82"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.94%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.529"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.35%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "522.89"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.56%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "7.06"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.92%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "102.57"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "13.34"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.46%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -6.05%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +4.56%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.81%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.19%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +11.43%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "32.80"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.84%  "
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "Bittensor"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "662.19"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +7.03%  "
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "Binance-PegBSC-USD"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.997"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.601"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.57%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.13"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +16.28%  "
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "41.54"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +25.63%  "
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.162"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +3.98%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.977"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +4.66%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.98"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.51%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.448"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.57%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0458"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.36%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.27%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.61"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.90%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.69%  "
